# report phan mo dau va chuong 1 OK
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 21-28 (A column checkboxes) as done (TRUE)
$ws.Range("A21:A28").Value = $true

# Move the selection / scroll position to row 34 area (B34), keeping the
# frozen header rows (1:3) in place, matching the saved view state
# (topLeftCell A37, active cell B34).
$ws.Range("B34").Select()
$excel.ActiveWindow.ScrollRow = 37

$wb.Save()
